# Fruta / hortaliza, semanal
# Inserts a new weekly price record for "Pepino dulce" (Vega Modelo de Temuco)
# at row 192, pushing the existing rows 192:271 down to 193:272.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 192 (shifts rows 192:271 down to 193:272,
# carrying their formatting/styles with them).
$ws.Rows.Item(192).Insert()

# Populate the newly inserted row 192 with the new weekly record.
$ws.Range("A192").Value = 10
$ws.Range("B192").Value = "Vega Modelo de Temuco"
$ws.Range("C192").Value = "La Araucanía"
$ws.Range("D192").Value = 44795
$ws.Range("E192").Value = 9
$ws.Range("F192").Value = 100112043
$ws.Range("G192").Value = "Pepino dulce"
$ws.Range("H192").Value = "Cultivar IV Región"
$ws.Range("I192").Value = "Primera"
$ws.Range("J192").Value = 300
$ws.Range("K192").Value = 18000
$ws.Range("L192").Value = 19000
$ws.Range("M192").Value = 18333
$ws.Range("N192").Value = "$/bandeja 18 kilos"
$ws.Range("O192").Value = "Provincia de Limarí"
$ws.Range("P192").Value = 1018
$ws.Range("Q192").Value = 18
$ws.Range("R192").Value = "Hortaliza"
